# Update gh-pages output figures (generated at 456a3b4)
# Updates the 'want to go' counts (column F) and refreshed cover image
# links (column I) across the four sheets of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 294
$ws.Range("F6").Value = 387
$ws.Range("F8").Value = 36
$ws.Range("F9").Value = 495
$ws.Range("F11").Value = 288
$ws.Range("F12").Value = 122
$ws.Range("F13").Value = 103
$ws.Range("F14").Value = 226
$ws.Range("F15").Value = 24
$ws.Range("F16").Value = 403
$ws.Range("F17").Value = 6554
$ws.Range("F19").Value = 68
$ws.Range("F21").Value = 7494
$ws.Range("F24").Value = 3367
$ws.Range("F25").Value = 18
$ws.Range("F26").Value = 1151
$ws.Range("F27").Value = 875
$ws.Range("F29").Value = 17
$ws.Range("F32").Value = 198
$ws.Range("F34").Value = 1578
$ws.Range("F35").Value = 6
$ws.Range("F36").Value = 144
$ws.Range("F39").Value = 1166
$ws.Range("F40").Value = 1679
$ws.Range("I15").Value = "//i1.hdslb.com/bfs/openplatform/202402/SjXuo1AJ1708338066210.png"

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 66
$ws.Range("F4").Value = 46
$ws.Range("F7").Value = 78

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1211

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1211
$ws.Range("F7").Value = 294
$ws.Range("F8").Value = 387
$ws.Range("F10").Value = 36
$ws.Range("F11").Value = 495
$ws.Range("F14").Value = 288
$ws.Range("F15").Value = 122
$ws.Range("F16").Value = 66
$ws.Range("F17").Value = 103
$ws.Range("F18").Value = 226
$ws.Range("F19").Value = 24
$ws.Range("F20").Value = 403
$ws.Range("F21").Value = 6554
$ws.Range("F23").Value = 68
$ws.Range("F25").Value = 7494
$ws.Range("F28").Value = 3367
$ws.Range("F29").Value = 18
$ws.Range("F30").Value = 1151
$ws.Range("F31").Value = 875
$ws.Range("F33").Value = 17
$ws.Range("F36").Value = 46
$ws.Range("F37").Value = 198
$ws.Range("F39").Value = 1578
$ws.Range("F40").Value = 6
$ws.Range("F41").Value = 144
$ws.Range("F44").Value = 1166
$ws.Range("F45").Value = 1679
$ws.Range("F49").Value = 78
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202402/SjXuo1AJ1708338066210.png"

Write-Host "Applied 58 cell updates across 4 sheets."
